$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.404345
$ws.Range("H2").Value = 40.213035
$ws.Range("I2").Value = 0.01122005832922476
$ws.Range("J2").Value = 0.01122005832922476
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 66.47763370022
$ws.Range("R2").Value = 598.2987033019799
$ws.Range("S2").Value = 0.00928085785146927
$ws.Range("T2").Value = 0.009280857851469272
$ws.Range("G3").Value = 13.404345
$ws.Range("H3").Value = 40.213035
$ws.Range("I3").Value = 0.01122005832922476
$ws.Range("J3").Value = 0.01122005832922476
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("Q3").Value = 7.57181065868
$ws.Range("R3").Value = 68.14629592812
$ws.Range("S3").Value = 0.001057090851313145
$ws.Range("T3").Value = 0.001057090851313145
$ws.Range("G4").Value = 13.404345
$ws.Range("H4").Value = 40.213035
$ws.Range("I4").Value = 0.01122005832922476
$ws.Range("J4").Value = 0.01122005832922476
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("Q4").Value = 6.31844184757
$ws.Range("R4").Value = 56.86597662813
$ws.Range("S4").Value = 0.0008821096264423439
$ws.Range("T4").Value = 0.000882109626442344
$ws.Range("I5").Value = 0.9315566574535661
$ws.Range("J5").Value = 0.9315566574535661
$ws.Range("M5").Value = 4.959409333333333
$ws.Range("N5").Value = 14.878228
$ws.Range("O5").Value = 0.8271666313262851
$ws.Range("P5").Value = 0.8271666313262852
$ws.Range("Q5").Value = 5519.372576156503
$ws.Range("R5").Value = 49674.35318540852
$ws.Range("S5").Value = 0.7705525822354403
$ws.Range("T5").Value = 0.7705525822354404
$ws.Range("I6").Value = 0.9315566574535661
$ws.Range("J6").Value = 0.9315566574535661
$ws.Range("O6").Value = 0.09421438109281059
$ws.Range("P6").Value = 0.09421438109281059
$ws.Range("S6").Value = 0.08776603393487509
$ws.Range("T6").Value = 0.08776603393487509
$ws.Range("I7").Value = 0.9315566574535661
$ws.Range("J7").Value = 0.9315566574535661
$ws.Range("O7").Value = 0.07861898758090437
$ws.Range("P7").Value = 0.07861898758090438
$ws.Range("S7").Value = 0.07323804128325069
$ws.Range("T7").Value = 0.0732380412832507
$ws.Range("I8").Value = 0.05722328421720919
$ws.Range("J8").Value = 0.05722328421720919
$ws.Range("M8").Value = 4.959409333333333
$ws.Range("N8").Value = 14.878228
$ws.Range("O8").Value = 0.8271666313262851
$ws.Range("P8").Value = 0.8271666313262852
$ws.Range("Q8").Value = 339.0417781881578
$ws.Range("R8").Value = 3051.37600369342
$ws.Range("S8").Value = 0.0473331912393755
$ws.Range("T8").Value = 0.04733319123937551
$ws.Range("I9").Value = 0.05722328421720919
$ws.Range("J9").Value = 0.05722328421720919
$ws.Range("O9").Value = 0.09421438109281059
$ws.Range("P9").Value = 0.09421438109281059
$ws.Range("S9").Value = 0.00539125630662236
$ws.Range("T9").Value = 0.00539125630662236
$ws.Range("I10").Value = 0.05722328421720919
$ws.Range("J10").Value = 0.05722328421720919
$ws.Range("O10").Value = 0.07861898758090437
$ws.Range("P10").Value = 0.07861898758090438
$ws.Range("S10").Value = 0.00449883667121133
$ws.Range("T10").Value = 0.00449883667121133
